$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.289457048391455
$ws.Range("C2").Value = 0.2136032845924944
$ws.Range("D2").Value = 0.07782535456573214
$ws.Range("E2").Value = 0.07189710192061671
$ws.Range("G2").Value = 1.545392280339001
$ws.Range("H2").Value = 1.346608427783679
$ws.Range("I2").Value = 1.260926226027763
$ws.Range("L2").Value = 0.2167030269477337
$ws.Range("M2").Value = 0.2805339692929323
$ws.Range("B3").Value = 1.195635153182423
$ws.Range("C3").Value = 0.186938350766269
$ws.Range("D3").Value = 0.07067906780041255
$ws.Range("E3").Value = 0.07184808947746646
$ws.Range("G3").Value = 1.520122562408119
$ws.Range("H3").Value = 1.341485981974557
$ws.Range("I3").Value = 1.258944554460683
$ws.Range("L3").Value = 0.2139232786851082
$ws.Range("M3").Value = 0.2660462923370659
$ws.Range("B4").Value = 1.138598068501608
$ws.Range("C4").Value = 0.1705102972586587
$ws.Range("D4").Value = 0.06633182752599964
$ws.Range("E4").Value = 0.07182008781319738
$ws.Range("G4").Value = 1.505706581488084
$ws.Range("H4").Value = 1.33910567828778
$ws.Range("I4").Value = 1.258477166854803
$ws.Range("L4").Value = 0.2123293458799722
$ws.Range("M4").Value = 0.2572953244346152
$ws.Range("B5").Value = 1.115498450131099
$ws.Range("C5").Value = 0.1638014751106311
$ws.Range("D5").Value = 0.06457037067971783
$ws.Range("E5").Value = 0.07180920367668531
$ws.Range("G5").Value = 1.500106900764507
$ws.Range("H5").Value = 1.338327434155985
$ws.Range("I5").Value = 1.258474422004475
$ws.Range("L5").Value = 0.2117081644317409
$ws.Range("M5").Value = 0.2537655810910024
$ws.Range("B6").Value = 1.111671445380409
$ws.Range("C6").Value = 0.1626866118544399
$ws.Range("D6").Value = 0.06427848700545269
$ws.Range("E6").Value = 0.07180742821163255
$ws.Range("G6").Value = 1.499193637958854
$ws.Range("H6").Value = 1.338209771260154
$ws.Range("I6").Value = 1.258485283135983
$ws.Range("L6").Value = 0.2116067300955962
$ws.Range("M6").Value = 0.2531816654446928
$ws.Range("B7").Value = 1.138285957914093
$ws.Range("C7").Value = 0.170419877793563
$ws.Range("D7").Value = 0.06630803126626006
$ws.Range("E7").Value = 0.07181993889171379
$ws.Range("G7").Value = 1.50562995108632
$ws.Range("H7").Value = 1.339094407027005
$ws.Range("I7").Value = 1.258476370725745
$ws.Range("L7").Value = 0.2123208536053909
$ws.Range("M7").Value = 0.2572475739032427
$ws.Range("B8").Value = 1.256989012852671
$ws.Range("C8").Value = 0.2044206073989585
$ws.Range("D8").Value = 0.07535280428193403
$ws.Range("E8").Value = 0.0718797685498167
$ws.Range("G8").Value = 1.536450065246299
$ws.Range("H8").Value = 1.344683092211369
$ws.Range("I8").Value = 1.26008698020722
$ws.Range("L8").Value = 0.2157211373879591
$ws.Range("M8").Value = 0.275508618933479
$ws.Range("B9").Value = 1.494295980468962
$ws.Range("C9").Value = 0.2706690688480933
$ws.Range("D9").Value = 0.09341872510248095
$ws.Range("E9").Value = 0.07201368130538555
$ws.Range("G9").Value = 1.605690676117149
$ws.Range("H9").Value = 1.361740584140364
$ws.Range("I9").Value = 1.269225220529236
$ws.Range("L9").Value = 0.2232858644346436
$ws.Range("M9").Value = 0.3124673864516794
$ws.Range("B10").Value = 1.67144045059041
$ws.Range("C10").Value = 0.3191070044301796
$ws.Range("D10").Value = 0.1069034015453099
$ws.Range("E10").Value = 0.07212218490540856
$ws.Range("G10").Value = 1.662041298714314
$ws.Range("H10").Value = 1.378033737584076
$ws.Range("I10").Value = 1.279634025894978
$ws.Range("L10").Value = 0.229393308477583
$ws.Range("M10").Value = 0.340327909052796
$ws.Range("B11").Value = 1.752643861237857
$ws.Range("C11").Value = 0.3410975400236111
$ws.Range("D11").Value = 0.1130863538956248
$ws.Range("E11").Value = 0.07217374760877249
$ws.Range("G11").Value = 1.68889145399416
$ws.Range("H11").Value = 1.386272233458044
$ws.Range("I11").Value = 1.285182447387314
$ws.Range("L11").Value = 0.232291805849357
$ws.Range("M11").Value = 0.3531575723979898
$ws.Range("B12").Value = 1.783482898304555
$ws.Range("C12").Value = 0.3494188199234145
$ws.Range("D12").Value = 0.1154348530633342
$ws.Range("E12").Value = 0.07219359015501237
$ws.Range("G12").Value = 1.699235655708094
$ws.Range("H12").Value = 1.389511540660862
$ws.Range("I12").Value = 1.287401300066648
$ws.Range("L12").Value = 0.2334067158074191
$ws.Range("M12").Value = 0.3580383061836585
$ws.Range("B13").Value = 1.776837198495684
$ws.Range("C13").Value = 0.3476269482370071
$ws.Range("D13").Value = 0.114928741446704
$ws.Range("E13").Value = 0.07218930261435563
$ws.Range("G13").Value = 1.696999962703302
$ws.Range("H13").Value = 1.388808568952555
$ws.Range("I13").Value = 1.28691817877889
$ws.Range("L13").Value = 0.233165829577004
$ws.Range("M13").Value = 0.3569861554021969
$ws.Range("B14").Value = 1.755179222305742
$ws.Range("C14").Value = 0.3417822564885
$ws.Range("D14").Value = 0.1132794222135374
$ws.Range("E14").Value = 0.0721753737162425
$ws.Range("G14").Value = 1.689738927545903
$ws.Range("H14").Value = 1.386536332894252
$ws.Range("I14").Value = 1.285362628417644
$ws.Range("L14").Value = 0.2323831830372427
$ws.Range("M14").Value = 0.353558663873315
$ws.Range("B15").Value = 1.741924685847437
$ws.Range("C15").Value = 0.3382014336327472
$ws.Range("D15").Value = 0.1122701019580177
$ws.Range("E15").Value = 0.07216688313033703
$ws.Range("G15").Value = 1.685314389250948
$ws.Range("H15").Value = 1.385160114669191
$ws.Range("I15").Value = 1.284425172505976
$ws.Range("L15").Value = 0.2319060448287331
$ws.Range("M15").Value = 0.3514621465706398
$ws.Range("B16").Value = 1.666146080634803
$ws.Range("C16").Value = 0.3176690076358
$ws.Range("D16").Value = 0.1065003235488575
$ws.Range("E16").Value = 0.07211885950265096
$ws.Range("G16").Value = 1.660311197908015
$ws.Range("H16").Value = 1.377512026262394
$ws.Range("I16").Value = 1.279287861197432
$ws.Range("L16").Value = 0.2292063058168026
$ws.Range("M16").Value = 0.3394925964115245
$ws.Range("B17").Value = 1.619817009088251
$ws.Range("C17").Value = 0.3050619295579224
$ws.Range("D17").Value = 0.1029733290261845
$ws.Range("E17").Value = 0.07208996296027337
$ws.Range("G17").Value = 1.645285249460045
$ws.Range("H17").Value = 1.373032391418946
$ws.Range("I17").Value = 1.276345200783254
$ws.Range("L17").Value = 0.2275809070082317
$ws.Range("M17").Value = 0.3321895738710055
$ws.Range("B18").Value = 1.593228062015783
$ws.Range("C18").Value = 0.2978065009430111
$ws.Range("D18").Value = 0.1009492706449606
$ws.Range("E18").Value = 0.07207354986773318
$ws.Range("G18").Value = 1.636757067837351
$ws.Range("H18").Value = 1.370533591046041
$ws.Range("I18").Value = 1.274729172494858
$ws.Range("L18").Value = 0.2266573317028389
$ws.Range("M18").Value = 0.3280037290083584
$ws.Range("B19").Value = 1.584235517110869
$ws.Range("C19").Value = 0.295349212545716
$ws.Range("D19").Value = 0.1002647407096049
$ws.Range("E19").Value = 0.07206802830699299
$ws.Range("G19").Value = 1.633889157285438
$ws.Range("H19").Value = 1.369700875662403
$ws.Range("I19").Value = 1.274195129220956
$ws.Range("L19").Value = 0.2263465668135893
$ws.Range("M19").Value = 0.3265889913282791
$ws.Range("B20").Value = 1.624742783246063
$ws.Range("C20").Value = 0.306404404134355
$ws.Range("D20").Value = 0.103348309006293
$ws.Range("E20").Value = 0.07209301757793796
$ws.Range("G20").Value = 1.646872940097381
$ws.Range("H20").Value = 1.373501202613284
$ws.Range("I20").Value = 1.276650527398033
$ws.Range("L20").Value = 0.2277527626271336
$ws.Range("M20").Value = 0.3329654756563585
$ws.Range("B21").Value = 1.761538277607144
$ws.Range("C21").Value = 0.3434991462462733
$ws.Range("D21").Value = 0.113763672055498
$ws.Range("E21").Value = 0.07217945637133028
$ws.Range("G21").Value = 1.691866861673731
$ws.Range("H21").Value = 1.387200493101233
$ws.Range("I21").Value = 1.28581632803531
$ws.Range("L21").Value = 0.232612595222335
$ws.Range("M21").Value = 0.354564792605828
$ws.Range("B22").Value = 1.851461211249216
$ws.Range("C22").Value = 0.3677075279857149
$ws.Range("D22").Value = 0.1206124542225382
$ws.Range("E22").Value = 0.07223779613758197
$ws.Range("G22").Value = 1.722303360986558
$ws.Range("H22").Value = 1.396851030146053
$ws.Range("I22").Value = 1.292493638716664
$ws.Range("L22").Value = 0.2358896980829854
$ws.Range("M22").Value = 0.3688118806351213
$ws.Range("B23").Value = 1.803420142231573
$ws.Range("C23").Value = 0.3547901812487453
$ws.Range("D23").Value = 0.1169532625274599
$ws.Range("E23").Value = 0.0722064900896271
$ws.Range("G23").Value = 1.705963960006159
$ws.Range("H23").Value = 1.39163632630715
$ws.Range("I23").Value = 1.288866706123699
$ws.Range("L23").Value = 0.2341314024778569
$ws.Range("M23").Value = 0.3611959763121959
$ws.Range("B24").Value = 1.622515696677283
$ws.Range("C24").Value = 0.3057974946164563
$ws.Range("D24").Value = 0.1031787691775463
$ws.Range("E24").Value = 0.07209163596286672
$ws.Range("G24").Value = 1.646154801298167
$ws.Range("H24").Value = 1.373289014511641
$ws.Range("I24").Value = 1.276512253350688
$ws.Range("L24").Value = 0.2276750327900317
$ws.Range("M24").Value = 0.3326146507154775
$ws.Range("B25").Value = 1.429609820719406
$ws.Range("C25").Value = 0.2527902012148218
$ws.Range("D25").Value = 0.08849504718347134
$ws.Range("E25").Value = 0.07197567837735352
$ws.Range("G25").Value = 1.586005214025107
$ws.Range("H25").Value = 1.356468766955004
$ws.Range("I25").Value = 1.26610780338283
$ws.Range("L25").Value = 0.2211430998897512
$ws.Range("M25").Value = 0.3023454069628286
